$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Move the "Assay" / "Sub-Aliquots" columns from D/E to J/K ---

# New header cells J1/K1, duplicating the look of the existing header row (bold font, same style as D1/E1)
$ws.Range("J1").Value = $ws.Range("D1").Value()
$ws.Range("K1").Value = $ws.Range("E1").Value()
$ws.Range("J1").Font.Bold = $true
$ws.Range("K1").Font.Bold = $true

# New data cells J2 (was D2), K2..K4 (were E2..E4)
$ws.Range("J2").Value = $ws.Range("D2").Value()
$ws.Range("K2").Value = $ws.Range("E2").Value()
$ws.Range("K3").Value = $ws.Range("E3").Value()
$ws.Range("K4").Value = $ws.Range("E4").Value()

# Clear the old D2 cell (column D is no longer used for data)
$ws.Range("D2").ClearContents()

# Column E now mirrors column F (Sample_id) instead of holding the old sub-aliquot ids
$ws.Range("E2").Value = $ws.Range("F2").Value()
$ws.Range("E3").Value = $ws.Range("F3").Value()
$ws.Range("E4").Value = $ws.Range("F4").Value()
